$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected, so it must be unprotected before any cell can be edited.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text (A9).
$old = $ws.Range("A9").Value()
$new = $old.Replace("2021-04-28", "2021-04-29")
$ws.Range("A9").Value = $new
# Editing this wrapped, multi-line cell can cause the engine to auto-expand the
# row height; auto-fit it back so row 9 keeps its original (default) height.
$ws.Rows(9).AutoFit()

# Update the Weight (D) and Percent Change (E) values for rows 2-6.
$ws.Range("D2").Value = 0.2515914468296035
$ws.Range("E2").Value = 0.006234264476681739

$ws.Range("D3").Value = 0.2501902662121832
$ws.Range("E3").Value = 0.01780745687256546

$ws.Range("D4").Value = 0.2463921555363704
$ws.Range("E4").Value = 0.007066444204534283

$ws.Range("D5").Value = 0.251826131421843
$ws.Range("E5").Value = 0.02766233766233772

$ws.Range("E6").Value = 0.01473095589423989

# Restore sheet protection.
$ws.Protect()
